$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Widen column A slightly to fit new, longer article names
$ws.Columns.Item(1).ColumnWidth = 53.29

# New article rows
$newRows = @(
    @{ Name = "16GB DDR4 3600MHz G.Skill Ripjaws V Series - DDR4 (2x8GB)"; Sku = 106547; Gtin = 4713294225634 },
    @{ Name = "16GB DDR4 3600MHz G.Skill Ripjaws V Series - DDR4 (1x16GB)"; Sku = 106538; Gtin = 4713294230089 },
    @{ Name = "32GB DDR5-6000 CL30 Kingston FURY Beast Kit 2x 16GB (AMD)"; Sku = 106525; Gtin = 740617343229 },
    @{ Name = "32GB DDR5-6000 CL30 Kingston FURY Beast Kit 2x 16GB (INTEL)"; Sku = 106524; Gtin = 740617342994 }
)

$row = 6
foreach ($item in $newRows) {
    $ws.Range("A5:B5").Copy() | Out-Null
    $ws.Range("A$row`:B$row").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($row, 1).Value = $item.Name
    $ws.Cells.Item($row, 2).Value = $item.Sku
    $ws.Cells.Item($row, 3).Value = $item.Gtin

    $row++
}

$ws.Range("B15").Select()
